# Updates cryptos list values (prices/volumes) and fixes the Polkadot/WrappedEther row order+data,
# matching the GitHub Actions scrape refresh described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.275.47"
$ws.Range("E2").Value = "  -5.03%  "

$ws.Range("D3").Value = "1.675.94"
$ws.Range("E3").Value = "  -1.88%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'218.37"
$ws.Range("E5").Value = "  -1.93%  "

$ws.Range("D6").Value = "'0.5126"
$ws.Range("E6").Value = "  -7.15%  "

$ws.Range("E7").Value = "  +0.13%  "

$ws.Range("E8").Value = "  +0.44%  "

$ws.Range("D9").Value = "'0.06397"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("D10").Value = "'21.51"
$ws.Range("E10").Value = "  -2.51%  "

$ws.Range("D11").Value = "'0.07381"
$ws.Range("E11").Value = "  -1.08%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.681.75"
$ws.Range("E12").Value = "  -2.24%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'4.567"
$ws.Range("E13").Value = "  +0.01%  "

$ws.Range("D14").Value = "'0.5822"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "1.903.60"
$ws.Range("E15").Value = "  -2.15%  "

$ws.Range("D16").Value = "'0.000008670"
$ws.Range("E16").Value = "  +3.50%  "

$ws.Range("D17").Value = "'64.90"
$ws.Range("E17").Value = "  -8.53%  "

$ws.Range("D18").Value = "26.345.92"
$ws.Range("E18").Value = "  -4.83%  "

$ws.Range("D19").Value = "'4.960"
$ws.Range("E19").Value = "  -4.03%  "

$ws.Range("E20").Value = "  +0.13%  "

$ws.Range("D21").Value = "'10.85"
$ws.Range("E21").Value = "  -0.91%  "

$ws.Range("D22").Value = "'190.52"
$ws.Range("E22").Value = "  -1.83%  "

$ws.Range("D23").Value = "'6.234"
$ws.Range("E23").Value = "  -2.50%  "

$ws.Range("E24").Value = "  -0.12%  "

$ws.Range("D25").Value = "'144.15"
$ws.Range("E25").Value = "  -2.81%  "

$ws.Range("D26").Value = "'7.672"
$ws.Range("E26").Value = "  -0.96%  "

$ws.Range("E27").Value = "  -0.08%  "

$ws.Range("D28").Value = "'15.67"
$ws.Range("E28").Value = "  +0.03%  "

$ws.Range("D29").Value = "'0.05925"
$ws.Range("E29").Value = "  -0.85%  "

$ws.Range("E30").Value = "  -5.99%  "

$ws.Range("D31").Value = "'1.326"
$ws.Range("E31").Value = "  -3.41%  "

$ws.Range("D32").Value = "'3.535"
$ws.Range("E32").Value = "  -1.42%  "

$ws.Range("D33").Value = "'3.520"
$ws.Range("E33").Value = "  -1.42%  "

$ws.Range("D34").Value = "'1.645"
$ws.Range("E34").Value = "  +0.33%  "

$ws.Range("E35").Value = "  +1.00%  "

$ws.Range("D36").Value = "'0.6028"
$ws.Range("E36").Value = "  -3.57%  "

$ws.Range("D37").Value = "'2.363"
$ws.Range("E37").Value = "  -2.47%  "

$ws.Range("E38").Value = "  -0.43%  "

$ws.Range("D39").Value = "'0.01621"
$ws.Range("E39").Value = "  -0.64%  "

$ws.Range("D40").Value = "'6.053"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("D41").Value = "1.080.53"
$ws.Range("E41").Value = "  -0.83%  "

$ws.Range("D42").Value = "'0.8691"
$ws.Range("E42").Value = "  -0.27%  "

$ws.Range("E43").Value = "  +0.55%  "

$ws.Range("D44").Value = "'99.88"
$ws.Range("E44").Value = "  +1.91%  "

$ws.Range("D45").Value = "1.824.39"
$ws.Range("E45").Value = "  -2.10%  "

$ws.Range("E46").Value = "  +1.70%  "

$ws.Range("D47").Value = "'56.08"
$ws.Range("E47").Value = "  -2.69%  "

$ws.Range("D48").Value = "'1.007"
$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("D49").Value = "'8.088"
$ws.Range("E49").Value = "  +1.88%  "

$ws.Range("D50").Value = "'0.4298"
$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("D51").Value = "'0.05200"
$ws.Range("E51").Value = "  -2.88%  "
